# "Generate Report for Archive"
# The localization status report is regenerated: the zh-cn / de-de entries
# on the Overview sheet have moved out of "Ready for handoff" and are now
# "In Translation", and the (now shorter) status text lets the Status
# columns on every sheet shrink to fit their contents.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: update the per-locale status values -------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

# --- zh-cn / de-de detail sheets: same status, it's the same underlying
#     shared-string value, just surfaced again on each locale's own sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- Narrow the Status columns to fit the new, shorter status text --------
# (was sized for "Ready for handoff"; now sized for "In Translation")
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
